$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 01:20"

# --- Update numeric data for various countries (case counts refreshed) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 161647
$ws.Range("C4").Value = 18156
$ws.Range("E4").Value = 153395
$ws.Range("F4").Value = 3512
$ws.Range("G4").Value = 415
$ws.Range("H4").Value = 2998

# Austria (row 17)
$ws.Range("B17").Value = 9618
$ws.Range("C17").Value = 830
$ws.Range("E17").Value = 8874

# Canada (row 18)
$ws.Range("B18").Value = 7428
$ws.Range("C18").Value = 1108
$ws.Range("E18").Value = 6249

# Australia (row 23)
$ws.Range("B23").Value = 4364
$ws.Range("C23").Value = 201
$ws.Range("E23").Value = 4101
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 19

# Ecuador (row 33)
$ws.Range("B33").Value = 1966
$ws.Range("C33").Value = 42
$ws.Range("E33").Value = 1901
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 62

# Barein (row 67)
$ws.Range("D67").Value = 295
$ws.Range("E67").Value = 216

# Tunez (row 78)
$ws.Range("D78").Value = 3
$ws.Range("E78").Value = 300
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 9

# Estado de Palestina (row 110)
$ws.Range("B110").Value = 117
$ws.Range("C110").Value = 8
$ws.Range("E110").Value = 98

# Trinidad yTobago (row 118)
$ws.Range("B118").Value = 83
$ws.Range("C118").Value = 5
$ws.Range("E118").Value = 79

# Angola (row 181)
$ws.Range("D181").Value = 1
$ws.Range("E181").Value = 4

# --- Reorder countries Sudan..Butan block (rows 187-194) and refresh Mauritania's data ---
# New order: Sudan, Mauritania, Montserrat, Islas Turcas y Caicos, Fiyi, Republica del Chad, Nepal, Butan

# Row 187: Sudan (unchanged values)
$ws.Range("A187").Value = "Sudan"
$ws.Range("B187").Value = 6
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 4
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 1
$ws.Range("H187").Value = 2

# Row 188: Mauritania (moved here, with updated case counts)
$ws.Range("A188").Value = "Mauritania"
$ws.Range("B188").Value = 6
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 2
$ws.Range("E188").Value = 3
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 1
$ws.Range("H188").Value = 1

# Row 189: Montserrat (moved here, values unchanged)
$ws.Range("A189").Value = "Montserrat"
$ws.Range("B189").Value = 5
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 5
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

# Row 190: Islas Turcas y Caicos (moved here, values unchanged)
$ws.Range("A190").Value = "Islas Turcas y Caicos"
$ws.Range("B190").Value = 5
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 5
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 191: Fiyi (moved here, values unchanged)
$ws.Range("A191").Value = "Fiyi"
$ws.Range("B191").Value = 5
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

# Row 192: Republica del Chad (moved here, values unchanged)
$ws.Range("A192").Value = "Republica del Chad"
$ws.Range("B192").Value = 5
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 0
$ws.Range("E192").Value = 5
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

# Row 193: Nepal (moved here, values unchanged)
$ws.Range("A193").Value = "Nepal"
$ws.Range("B193").Value = 5
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 1
$ws.Range("E193").Value = 4
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

# Row 194: Butan (unchanged)
$ws.Range("A194").Value = "Butan"
$ws.Range("B194").Value = 4
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 4
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0
